$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Prepare rows 23 and 24 by duplicating formatting from row 22 ---
$ws.Range("A22:X22").Copy()
$ws.Range("A23:X23").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("A22:X22").Copy()
$ws.Range("A24:X24").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# --- Populate new text values in the precise order the strings were
#     originally authored, so the shared-string table ends up in the
#     same order as the source workbook. ---
$ws.Range("F23").Value = 'Trips lab notebook Will driver be late? PREP V11 MODEL V20'
$ws.Range("G23").Value = 'test for leakage train on older data and test on newer data'
$ws.Range("X22").Value = 'all scores improved. medabserror down 17 seconds to 3min42sec. Explaining 26% of the variance up from 20%. New feature was #1 kicking driver_id down to #7. Top 10 features used: [''avg_prior_arrived_late_seconds'', ''claimed_before_trip_start_secs'', ''origin_location_id'', ''driver_previous_completed_trips'', ''driver_home_lon'', ''driver_home_lat'', ''driver_id'', ''scheduled_starts_at_pdt_hour'', ''driver_created_at_pdt_dayofyear'', ''destination_location_id'']  Note that 28 features were above 1% importance, indicating a long tail of causes. Now 5 of top 10 are about driver (independent of the trip) and 5 are about trip. Precision of the late classes in 40-50% range!'
$ws.Range("X23").Value = 'Model did not perform as well on this test set over a seprate, later time frame. Explaining 19% of variance compared to 26% in random test set. Precision of late classes in 20%-30% range. Leakage existed in some way in randomly selected test set. Possibly due to different drivers driving in different times of the year. Top 10 features used: [''avg_prior_arrived_late_seconds'', ''claimed_before_trip_start_secs'', ''driver_previous_completed_trips'', ''driver_home_lon'', ''driver_id'', ''driver_home_lat'', ''origin_location_id'', ''scheduled_starts_at_pdt_hour'', ''driver_created_at_pdt_dayofyear'', ''origin_lon'']. Similar importance features.'
$ws.Range("B24").Value = 'Trips lab notebook Will driver be late? PREP V12'
$ws.Range("C24").Value = 'added feature avg_prior_arrived_late_seconds_to_origin_location'
$ws.Range("F24").Value = 'Trips lab notebook Will driver be late? PREP V12 MODEL V21'
$ws.Range("G24").Value = 'continute to use time series test set and test new feature'
$ws.Range("X24").Value = 'Model performance improved not back up to level with randomly selected test set. Explained 19% of variance. New feature was used #6 rank importance. Interesting that origin_location which was #3 now #8 and driver_id is #3, up from #5. This indicates there are still driver and location patterns unrelated to past performance of target variable (lateness). Precision of late classes in high 20%s. The average precision/recall/F1 are dominated by the large early arrival class. Top 10 features used: [''avg_prior_arrived_late_seconds'', ''claimed_before_trip_start_secs'', ''driver_id'', ''driver_previous_completed_trips'', ''driver_home_lat'', ''avg_prior_arrived_late_seconds_to_origin_location'', ''driver_home_lon'', ''origin_location_id'', ''scheduled_starts_at_pdt_hour'', ''destination_location_id'']'

# --- Fill in the remaining (reused / numeric) values for row 23 ---
$ws.Range("A23").Value = 43415
$ws.Range("B23").Value = 'Trips lab notebook Will driver be late? PREP V11'
$ws.Range("C23").Value = 'same'
$ws.Range("D23").Value = 125675
$ws.Range("E23").Value = 215
$ws.Range("H23").Value = 'GradientBoostingRegressor'
$ws.Range("I23").Value = 'max_depth=4, random_state=808, n_estimators=201'
$ws.Range("J23").Value = 100198
$ws.Range("K23").Value = 25477
$ws.Range("L23").Value = -0.0086116083461328
$ws.Range("M23").Value = 0.18
$ws.Range("N23").Value = 7.57
$ws.Range("O23").Value = 206403
$ws.Range("P23").Value = 0.17899999999999999
$ws.Range("Q23").Value = 337
$ws.Range("R23").Value = 251
$ws.Range("S23").Value = 4.18
$ws.Range("T23").Value = 0.75
$ws.Range("U23").Value = 0.83
$ws.Range("V23").Value = 0.77
$ws.Range("W23").Value = 152
$ws.Rows("23").RowHeight = 135

# --- Fill in the remaining (reused / numeric) values for row 24 ---
$ws.Range("A24").Value = 43415
$ws.Range("D24").Value = 125675
$ws.Range("E24").Value = 216
$ws.Range("H24").Value = 'GradientBoostingRegressor'
$ws.Range("I24").Value = 'max_depth=4, random_state=808, n_estimators=201'
$ws.Range("J24").Value = 100198
$ws.Range("K24").Value = 25477
$ws.Range("L24").Value = -0.00831184410557102
$ws.Range("M24").Value = 0.186
$ws.Range("N24").Value = 7.54
$ws.Range("O24").Value = 204536
$ws.Range("P24").Value = 0.186
$ws.Range("Q24").Value = 334
$ws.Range("R24").Value = 245
$ws.Range("S24").Value = 4.09
$ws.Range("T24").Value = 0.75
$ws.Range("U24").Value = 0.83
$ws.Range("V24").Value = 0.77
$ws.Range("W24").Value = 145
$ws.Rows("24").RowHeight = 165

# --- Update selection to mirror where the author ended up after editing ---
$ws.Range("X25").Select()

Write-Host "Edit complete"
